$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '40.051.40'
$ws.Cells.Item(2, 5).Value = '  +2.59%  '

$ws.Cells.Item(3, 4).Value = '2.237.07'
$ws.Cells.Item(3, 5).Value = '  +1.03%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '294.62'
$ws.Cells.Item(5, 5).Value = '  -0.51%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '86.06'
$ws.Cells.Item(6, 5).Value = '  +7.05%  '

$ws.Cells.Item(7, 5).Value = '  +1.97%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.472'
$ws.Cells.Item(9, 5).Value = '  +3.14%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '31.06'
$ws.Cells.Item(10, 5).Value = '  +10.85%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0789'
$ws.Cells.Item(11, 5).Value = '  +2.07%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '46.94'
$ws.Cells.Item(12, 5).Value = '  +2.42%  '

$ws.Cells.Item(13, 5).Value = '  +1.20%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.45'
$ws.Cells.Item(14, 5).Value = '  +5.62%  '

$ws.Cells.Item(15, 4).Value = '2.582.93'
$ws.Cells.Item(15, 5).Value = '  +1.23%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.12'
$ws.Cells.Item(16, 5).Value = '  +1.47%  '

$ws.Cells.Item(17, 4).Value = '2.291.18'
$ws.Cells.Item(17, 5).Value = '  +2.68%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.729'
$ws.Cells.Item(18, 5).Value = '  +2.62%  '

$ws.Cells.Item(19, 4).Value = '39.986.44'
$ws.Cells.Item(19, 5).Value = '  +2.70%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0890'
$ws.Cells.Item(20, 5).Value = '  +3.69%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.80'
$ws.Cells.Item(21, 5).Value = '  +1.96%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.81'
$ws.Cells.Item(22, 5).Value = '  +10.47%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '65.34'
$ws.Cells.Item(23, 5).Value = '  +1.00%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '234.90'
$ws.Cells.Item(24, 5).Value = '  +4.46%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.999'
$ws.Cells.Item(25, 5).Value = '  -0.08%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.47'
$ws.Cells.Item(26, 5).Value = '  +3.65%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.84'
$ws.Cells.Item(27, 5).Value = '  +5.27%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '22.80'
$ws.Cells.Item(28, 5).Value = '  +2.41%  '

$ws.Cells.Item(29, 5).Value = '  +3.09%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.22'
$ws.Cells.Item(30, 5).Value = '  +3.72%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '33.25'
$ws.Cells.Item(31, 5).Value = '  +6.82%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '152.67'
$ws.Cells.Item(32, 5).Value = '  +2.82%  '

$ws.Cells.Item(33, 5).Value = '  +0.08%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.87'
$ws.Cells.Item(34, 5).Value = '  +2.79%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0719'
$ws.Cells.Item(35, 5).Value = '  +4.78%  '

$ws.Cells.Item(36, 5).Value = '  +2.31%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '16.32'
$ws.Cells.Item(37, 5).Value = '  +13.99%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.111'
$ws.Cells.Item(38, 5).Value = '  +2.79%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.1000'
$ws.Cells.Item(39, 5).Value = '  +4.33%  '

$ws.Cells.Item(40, 5).Value = '  +3.48%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.69'
$ws.Cells.Item(41, 5).Value = '  +6.23%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.83'
$ws.Cells.Item(42, 5).Value = '  +6.05%  '

$ws.Cells.Item(43, 4).Value = '2.042.54'
$ws.Cells.Item(43, 5).Value = '  +7.70%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.23'
$ws.Cells.Item(44, 5).Value = '  +7.60%  '

$ws.Cells.Item(45, 5).Value = '  +6.09%  '

$ws.Cells.Item(46, 5).Value = '  +13.82%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '16.32'
$ws.Cells.Item(47, 5).Value = '  +1.65%  '

$ws.Cells.Item(48, 5).Value = '  +2.41%  '

$ws.Cells.Item(49, 4).Value = '2.452.92'
$ws.Cells.Item(49, 5).Value = '  +1.31%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '70.77'
$ws.Cells.Item(50, 5).Value = '  +1.40%  '

$ws.Cells.Item(51, 5).Value = '  +14.52%  '
